# Update the "All-time" (B) stat values on the Pogdude69 and one4kat sheets
# with a newly-uploaded snapshot, and convert the Session/Daily/Yesterday/
# Monthly "Delta" columns (C, E, G, I) from live formulas (=B-D, =B-F, =B-H,
# =B-J) into their static computed results, matching the uploaded snapshot.
#
# D, F, H, J (the Snapshot columns) are left untouched.
#
# Each hashtable below maps row number -> @(NewB, NewC, NewE, NewG, NewI).

$pogdudeData = @{
     2 = @(109, 62, 57, 63, 63)
     3 = @(1577695.833333331, 77141.24999999977, 75811.24999999977, 77211.66666666674, 77211.66666666674)
     4 = @(324, 16, 16, 16, 16)
     5 = @(158082, 4619, 4619, 4619, 4619)
     6 = @(591157, 30886, 30413, 30940, 30940)
     7 = @(8746, 222, 219, 223, 223)
     8 = @(2186, 55, 54, 56, 56)
     9 = @(18418, 960, 943, 961, 961)
    10 = @(6113, 131, 130, 132, 132)
    11 = @(228728, 10603, 10426, 10623, 10623)
    12 = @(153, 1, 1, 1, 1)
    13 = @(6222, 164, 162, 164, 164)
    14 = @(12281, 824, 808, 824, 824)
    15 = @(27473, 1552, 1525, 1553, 1553)
    16 = @(18980, 1012, 994, 1013, 1013)
    17 = @(1923, 104, 103, 104, 104)
    18 = @(7287, 446, 439, 446, 446)
    19 = @(105, 8, 8, 8, 8)
    20 = @(61, 1, 1, 1, 1)
    21 = @(964, 86, 84, 86, 86)
    22 = @(1837284, 93363, 91791, 93532, 93532)
}

$one4katData = @{
     2 = @(49, 39, 39, 39, 39)
     3 = @(588488.7499999991, 15374.58333333326, 15374.58333333326, 15374.58333333326, 15374.58333333326)
     4 = @(122, 3, 3, 3, 3)
     5 = @(141282, 2282, 2282, 2282, 2282)
     6 = @(100793, 2568, 2568, 2568, 2568)
     7 = @(3530, 90, 90, 90, 90)
     8 = @(818, 20, 20, 20, 20)
     9 = @(5932, 158, 158, 158, 158)
    10 = @(2023, 24, 24, 24, 24)
    11 = @(56330, 1083, 1083, 1083, 1083)
    12 = @(53, 1, 1, 1, 1)
    13 = @(2585, 66, 66, 66, 66)
    14 = @(3904, 133, 133, 133, 133)
    15 = @(4207, 108, 108, 108, 108)
    16 = @(2759, 62, 62, 62, 62)
    17 = @(1526, 4, 4, 4, 4)
    18 = @(1356, 44, 44, 44, 44)
    19 = @(11, 0, 0, 0, 0)
    20 = @(9, 0, 0, 0, 0)
    21 = @(81, 2, 2, 2, 2)
    22 = @(622398, 15239, 15239, 15239, 15239)
}

function Update-StatsSheet {
    param($SheetName, $RowData)

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in $RowData.Keys) {
        $vals = $RowData[$row]

        $ws.Cells.Item($row, 2).Value = $vals[0]   # B - All-time
        $ws.Cells.Item($row, 3).Value = $vals[1]   # C - Session Delta
        $ws.Cells.Item($row, 5).Value = $vals[2]   # E - Daily Delta
        $ws.Cells.Item($row, 7).Value = $vals[3]   # G - Yesterday Delta
        $ws.Cells.Item($row, 9).Value = $vals[4]   # I - Monthly Delta
    }
}

$wb = $excel.ActiveWorkbook

Update-StatsSheet "Pogdude69" $pogdudeData
Update-StatsSheet "one4kat" $one4katData
